$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header (H1) - copy formatting from the neighboring
# header cell (G1) so it picks up the same bold/border/center-top style,
# then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data values for the Save column.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
